# [FIX] budget overview report
# Insert a new "Budget Method" label row into the report's parameter list
# (right above "Activity Group"), matching the rest of the label rows above it.
# All subsequent rows - including the blank spacer row and the column-header
# table row further down the sheet - shift down by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above row 11 ("Activity Group"); Excel carries the
# formatting of the row above down into the newly inserted row automatically.
$ws.Rows(11).Insert()

# Populate the new row with its label.
$ws.Range("A11").Value = "Budget Method"
